# add nowcasts for 2025q4
# Shift the date labels forward one quarter and refresh the nowcast /
# revision figures for rows 2-7 with the latest 2025Q4 estimates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the vintage dates as plain text (not real dates).
# Force a text format first so Excel does not auto-convert the
# "yyyy-mm-dd" strings into date serial numbers, then restore the
# default "Normal" style so no extra formatting is left behind.
$ws.Range("A2:A7").NumberFormat = "@"
$ws.Range("A2").Value = "2025-09-30"
$ws.Range("A3").Value = "2025-10-15"
$ws.Range("A4").Value = "2025-10-30"
$ws.Range("A5").Value = "2025-11-15"
$ws.Range("A6").Value = "2025-11-30"
$ws.Range("A7").Value = "2025-12-15"
$ws.Range("A2:A7").Style = "Normal"

# --- Row 2 nowcast values ---
$ws.Range("B2").Value = 0.29024304835761255
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# --- Row 3 nowcast values ---
$ws.Range("B3").Value = 0.28930541399418097
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.005245537116881138
$ws.Range("E3").Value = -0.00014769779747891948
$ws.Range("F3").Value = -0.0003671430698622386
$ws.Range("G3").Value = -0.00025345890231253645
$ws.Range("H3").Value = 0.00003708225738461833
$ws.Range("I3").Value = 0.00000007764127736660174
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -0.000002194351012174689

# --- Row 4 nowcast values ---
$ws.Range("B4").Value = 0.29867716825012003
$ws.Range("C4").Value = 0.006598900314074544
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -0.000004663855265913454
$ws.Range("F4").Value = -0.0000013198676655936466
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.00020968419938015891
$ws.Range("I4").Value = -0.0006584834791823304
$ws.Range("J4").Value = 0.0028526301968664213
$ws.Range("K4").Value = 0.000000386663450480107

# --- Row 5 nowcast values ---
$ws.Range("B5").Value = 0.29786666036096704
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = -0.002936113887931363
$ws.Range("E5").Value = -0.0012264259332025148
$ws.Range("F5").Value = 0.0014323113802256445
$ws.Range("G5").Value = 0.00007472552015529076
$ws.Range("H5").Value = 0.00018578385978718513
$ws.Range("I5").Value = -0.00007524904150874597
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.00040081536563496556

# --- Row 6 nowcast values ---
$ws.Range("B6").Value = 0.3203875304584916
$ws.Range("C6").Value = -0.009652584648399735
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.000006266343421440825
$ws.Range("F6").Value = 0.00020762867863823926
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.0001909689549819912
$ws.Range("I6").Value = 0.001429306604338888
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -0.0000028783939146070914

# --- Row 7 nowcast values ---
$ws.Range("B7").Value = 0.32355527914012416
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.005703312244309744
$ws.Range("E7").Value = 0.00025610445746612206
$ws.Range("F7").Value = 0.001321695229232745
$ws.Range("G7").Value = 0.0028286664371489586
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.00009090804944145159

# Columns C, D, H and J became slightly narrower in this update.
$ws.Columns.Item(3).ColumnWidth = 14.333333333333334
$ws.Columns.Item(4).ColumnWidth = 14.333333333333334
$ws.Columns.Item(8).ColumnWidth = 14.833333333333334
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
